# Fill in the "Decision" column (column 3) for rows whose decision cell
# is currently empty, using the LPA reference in column 1 to pick the
# right row so the mapping is unambiguous.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$decisions = @{
    "3/22/0005/LBC"  = "Grant Consent subject to Conditions"
    "3/21/2884/HH"   = "Grant Planning Permission Subject to Conditions"
    "3/21/3023/HH"   = "Grant Planning Permission Subject to Conditions"
    "3/22/0316/HH"   = "Grant Planning Permission Subject to Conditions"
    "3/22/0340/HH"   = "Grant Planning Permission Subject to Conditions"
    "3/21/1423/VAR"  = "Refuse Planning Permission"
    "3/22/0291/HH"   = "Grant Planning Permission Subject to Conditions"
    "3/22/0286/HH"   = "Grant Planning Permission Subject to Conditions"
    "3/22/0115/HH"   = "Grant Planning Permission Subject to Conditions"
    "3/21/2573/FUL"  = "Grant Planning Permission Subject to Conditions"
    "3/22/0273/ARPN" = "Prior Approval is Required and Granted Subject to Conditions"
    "3/22/0267/HH"   = "Grant Planning Permission Subject to Conditions"
    "3/21/1405/VAR"  = "Grant Planning Permission Subject to Conditions"
    "3/21/2620/LBC"  = "Grant Consent subject to Conditions"
}

for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $refCell = $t.Cell($r, 1).Range.Text
    $ref = $refCell.TrimEnd([char]13, [char]7)

    if ($decisions.ContainsKey($ref)) {
        $decisionCell = $t.Cell($r, 3)
        $cellRange = $decisionCell.Range
        $cellRange.End = $cellRange.End - 1
        $cellRange.Text = $decisions[$ref]
    }
}
